$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "parameters": add a new "cavalry" parameter block (distance +
# height-difference threshold) right before the existing "flier"
# parameters, pushing flier_distance / flier_distance_height_gain down
# two rows.
# ------------------------------------------------------------------
$wsParams = $wb.Worksheets.Item("parameters")
$wsParams.Rows("8:9").Insert(-4121)   # xlShiftDown

$wsParams.Range("A8").Value = "cavalry_distance"
$wsParams.Range("B8").Value = 4.5
$wsParams.Range("A9").Value = "cavalry_height_difference_threshold"
$wsParams.Range("B9").Value = 2

# ------------------------------------------------------------------
# Sheet "interactions": add a new "cavalry" column right before the
# existing "siege" column, with values for the b1/b2/b3 and h1/h2 rows.
# ------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("interactions")
$wsInt.Columns("F").Insert(-4161)   # xlShiftToRight

$wsInt.Range("F1").Value = "cavalry"
$wsInt.Range("F2").Value = 2
$wsInt.Range("F3").Value = 2
$wsInt.Range("F4").Value = 2
$wsInt.Range("F5").Value = -1
$wsInt.Range("F6").Value = -1

# ------------------------------------------------------------------
# Restore/update sheet selections and the active tab: "parameters"
# keeps a stored selection over the new cavalry rows, while
# "interactions" becomes the active sheet with F1 selected.
# ------------------------------------------------------------------
$wsParams.Activate()
$wsParams.Range("A8:B9").Select()

$wsInt.Activate()
$wsInt.Range("F1").Select()
